$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 11
$ws.Range("H11").Value = 79.42856999999999
$ws.Range("I11").Value = 79.42856999999999
$ws.Range("K11").Value = 79.42856999999999
$ws.Range("M11").Value = 60.57143000000001

# row 64
$ws.Range("H64").Value = 2847
$ws.Range("I64").Value = 2996
$ws.Range("K64").Value = 2996
$ws.Range("M64").Value = -2748

# row 67
$ws.Range("H67").Value = 2847
$ws.Range("I67").Value = 2996
$ws.Range("K67").Value = 2996
$ws.Range("M67").Value = -2138

# row 76
$ws.Range("H76").Value = 4750
$ws.Range("I76").Value = 3500
$ws.Range("K76").Value = 3500
$ws.Range("M76").Value = -3185

# row 79
$ws.Range("H79").Value = 4750
$ws.Range("I79").Value = 3500
$ws.Range("K79").Value = 3500
$ws.Range("M79").Value = -2408

# row 106
$ws.Range("H106").Value = 16449.5
$ws.Range("I106").Value = 29999
$ws.Range("K106").Value = 29999
$ws.Range("M106").Value = -29368

# row 132
$ws.Range("H132").Value = 100004024
$ws.Range("I132").Value = 100004024
$ws.Range("K132").Value = 300012072
$ws.Range("M132").Value = -300009542

# row 138
$ws.Range("H138").Value = 1563.6765
$ws.Range("I138").Value = 1694.8334
$ws.Range("K138").Value = 5084.5002
$ws.Range("M138").Value = 55.4997999999996

$ws = $wb.Worksheets.Item("ARM")
# row 45
$ws.Range("H45").Value = 2524.2856
$ws.Range("I45").Value = 1631.8182
$ws.Range("K45").Value = 1631.8182
$ws.Range("M45").Value = -1254.8182

# row 61
$ws.Range("H61").Value = 1874.5
$ws.Range("I61").Value = 1732.6666
$ws.Range("K61").Value = 1732.6666
$ws.Range("M61").Value = -1520.6666

# row 74
$ws.Range("H74").Value = 3817.2778
$ws.Range("I74").Value = 3478.2727
$ws.Range("J74").Value = 4350
$ws.Range("K74").Value = 3478.2727
$ws.Range("L74").Value = 4350
$ws.Range("M74").Value = -2604.2727
$ws.Range("N74").Value = -6098

# row 77
$ws.Range("H77").Value = 3817.2778
$ws.Range("I77").Value = 3478.2727
$ws.Range("J77").Value = 4350
$ws.Range("K77").Value = 17391.3635
$ws.Range("L77").Value = 21750
$ws.Range("M77").Value = -13023.3635
$ws.Range("N77").Value = -30486

# row 136
$ws.Range("H136").Value = 1874.5
$ws.Range("I136").Value = 1732.6666
$ws.Range("K136").Value = 5197.9998
$ws.Range("M136").Value = -2647.9998

$ws = $wb.Worksheets.Item("BSM")
# row 20
$ws.Range("H20").Value = 1949.4375
$ws.Range("I20").Value = 2027.9286
$ws.Range("J20").Value = 1400
$ws.Range("K20").Value = 2027.9286
$ws.Range("L20").Value = 1400
$ws.Range("M20").Value = -1780.9286
$ws.Range("N20").Value = -1894

# row 64
$ws.Range("H64").Value = 4749
$ws.Range("J64").Value = 5583.3335
$ws.Range("L64").Value = 5583.3335
$ws.Range("N64").Value = -6033.3335

# row 67
$ws.Range("H67").Value = 4749
$ws.Range("J67").Value = 5583.3335
$ws.Range("L67").Value = 5583.3335
$ws.Range("N67").Value = -7143.3335

# row 105
$ws.Range("H105").Value = 3087.2
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

# row 134
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 3336.7778
$ws.Range("I31").Value = 4299.6665
$ws.Range("J31").Value = 2855.3333
$ws.Range("K31").Value = 4299.6665
$ws.Range("L31").Value = 2855.3333
$ws.Range("M31").Value = -4004.6665
$ws.Range("N31").Value = -3445.3333

# row 34
$ws.Range("H34").Value = 3336.7778
$ws.Range("I34").Value = 4299.6665
$ws.Range("J34").Value = 2855.3333
$ws.Range("K34").Value = 4299.6665
$ws.Range("L34").Value = 2855.3333
$ws.Range("M34").Value = -4097.6665
$ws.Range("N34").Value = -3259.3333

# row 86
$ws.Range("H86").Value = 9945
$ws.Range("I86").Value = 4424.625
$ws.Range("J86").Value = 24666
$ws.Range("K86").Value = 4424.625
$ws.Range("L86").Value = 24666
$ws.Range("M86").Value = -3301.625
$ws.Range("N86").Value = -26912

# row 89
$ws.Range("H89").Value = 9945
$ws.Range("I89").Value = 4424.625
$ws.Range("J89").Value = 24666
$ws.Range("K89").Value = 22123.125
$ws.Range("L89").Value = 123330
$ws.Range("M89").Value = -16507.125
$ws.Range("N89").Value = -134562

# row 99
$ws.Range("H99").Value = 11506
$ws.Range("I99").Value = 11506
$ws.Range("K99").Value = 11506
$ws.Range("M99").Value = -10008

# row 126
$ws.Range("H126").Value = 11506
$ws.Range("I126").Value = 11506
$ws.Range("K126").Value = 34518
$ws.Range("M126").Value = -32048

# row 132
$ws.Range("H132").Value = 4547.625
$ws.Range("I132").Value = 4321.625
$ws.Range("K132").Value = 12964.875
$ws.Range("M132").Value = -10434.875

$ws = $wb.Worksheets.Item("CUL")
# row 11
$ws.Range("H11").Value = 471.25
$ws.Range("I11").Value = 294.33334
$ws.Range("J11").Value = 1002
$ws.Range("K11").Value = 883.0000200000001
$ws.Range("L11").Value = 3006
$ws.Range("M11").Value = -743.0000200000001
$ws.Range("N11").Value = -3286

# row 22
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").ClearContents()
$ws.Range("N22").Value = 0

# row 27
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").ClearContents()
$ws.Range("N27").Value = 0

# row 75
$ws.Range("H75").Value = 856.3333
$ws.Range("J75").Value = 856.3333
$ws.Range("L75").Value = 2568.9999
$ws.Range("N75").Value = -4564.9999

# row 78
$ws.Range("H78").Value = 856.3333
$ws.Range("J78").Value = 856.3333
$ws.Range("L78").Value = 7706.9997
$ws.Range("N78").Value = -17690.9997

# row 113
$ws.Range("H113").Value = 1119.1296
$ws.Range("I113").Value = 1132.9387
$ws.Range("K113").Value = 3398.8161
$ws.Range("M113").Value = -1228.8161

# row 120
$ws.Range("H120").Value = 800
$ws.Range("I120").Value = 800
$ws.Range("K120").Value = 2400
$ws.Range("M120").Value = 2438

$ws = $wb.Worksheets.Item("GSM")
# row 2
$ws.Range("H2").Value = 11.125
$ws.Range("I2").Value = 9.25
$ws.Range("K2").Value = 9.25
$ws.Range("M2").Value = 103.75

# row 46
$ws.Range("H46").Value = 29299.625
$ws.Range("J46").Value = 30259.666
$ws.Range("L46").Value = 30259.666
$ws.Range("N46").Value = -30571.666

# row 122
$ws.Range("H122").Value = 2422.2856
$ws.Range("I122").Value = 2422.2856
$ws.Range("K122").Value = 7266.8568
$ws.Range("M122").Value = -4816.8568

# row 132
$ws.Range("H132").Value = 3131.5356
$ws.Range("I132").Value = 3111.5
$ws.Range("J132").Value = 3251.75
$ws.Range("K132").Value = 9334.5
$ws.Range("L132").Value = 9755.25
$ws.Range("M132").Value = -6804.5
$ws.Range("N132").Value = -14815.25

$ws = $wb.Worksheets.Item("LTW")
# row 16
$ws.Range("H16").Value = 462.92307
$ws.Range("I16").Value = 553.5
$ws.Range("J16").Value = 161
$ws.Range("K16").Value = 553.5
$ws.Range("L16").Value = 161
$ws.Range("M16").Value = -383.5
$ws.Range("N16").Value = -501

# row 22
$ws.Range("H22").Value = 3933
$ws.Range("I22").Value = 899.5
$ws.Range("K22").Value = 899.5
$ws.Range("M22").Value = -604.5

# row 27
$ws.Range("H27").Value = 3933
$ws.Range("I27").Value = 899.5
$ws.Range("K27").Value = 899.5
$ws.Range("M27").Value = -792.5

# row 46
$ws.Range("H46").Value = 1201.8276
$ws.Range("I46").Value = 1209.7858
$ws.Range("K46").Value = 1209.7858
$ws.Range("M46").Value = -1021.7858

# row 93
$ws.Range("H93").Value = 998.5
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()

# row 122
$ws.Range("H122").Value = 2671.5
$ws.Range("I122").Value = 2505.9
$ws.Range("K122").Value = 7517.700000000001
$ws.Range("M122").Value = -5067.700000000001

# row 132
$ws.Range("H132").Value = 3995.75
$ws.Range("I132").Value = 3995.75
$ws.Range("K132").Value = 11987.25
$ws.Range("M132").Value = -9457.25

$ws = $wb.Worksheets.Item("WVR")
# row 113
$ws.Range("H113").Value = 445
$ws.Range("I113").Value = 369.22223
$ws.Range("K113").Value = 1107.66669
$ws.Range("M113").Value = 1062.33331

# row 132
$ws.Range("H132").Value = 4082.6667
$ws.Range("I132").Value = 4119.2
$ws.Range("J132").Value = 3900
$ws.Range("K132").Value = 12357.6
$ws.Range("L132").Value = 11700
$ws.Range("M132").Value = -9827.599999999999
$ws.Range("N132").Value = -16760
